# Add season record columns (Wins / Losses / Ties) to the right of the
# existing "Unnamed: 28" column (AC) on the player stats sheet.
#
# New header cells: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties"
# For every data row (2 through 45) fill in the team's season record:
#   Wins = 91, Losses = 71, Ties = 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered look used by the rest of row 1
# by copying the formatting from the neighboring header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$lastRow = 45
$wins = 91
$losses = 71
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
